$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.069504995555228
$ws.Range("D2").Value = 1.068227240103944
$ws.Range("E2").Value = 1.073200654573546
$ws.Range("F2").Value = 1.081748225534817
$ws.Range("I2").Value = 1.047873759189906
$ws.Range("J2").Value = 1.074439382078949
$ws.Range("K2").Value = 1.070932987304475
$ws.Range("L2").Value = 1.075893145079647
$ws.Range("M2").Value = 1.084418239958538

$ws.Range("B3").Value = 1.019999999999999
$ws.Range("C3").Value = 1.071079478604033
$ws.Range("D3").Value = 1.069436886922873
$ws.Range("E3").Value = 1.074584481136351
$ws.Range("F3").Value = 1.083168121359609
$ws.Range("I3").Value = 1.048271202457672
$ws.Range("J3").Value = 1.075668105253649
$ws.Range("K3").Value = 1.071957593096609
$ws.Range("L3").Value = 1.077092466562213
$ws.Range("M3").Value = 1.0856551839553

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.072096838609928
$ws.Range("D4").Value = 1.070218086496218
$ws.Range("E4").Value = 1.075478763245119
$ws.Range("F4").Value = 1.084085738166866
$ws.Range("I4").Value = 1.04852603424059
$ws.Range("J4").Value = 1.076461324540364
$ws.Range("K4").Value = 1.072618489562521
$ws.Range("L4").Value = 1.077866825694775
$ws.Range("M4").Value = 1.086453886992458

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.072524201142357
$ws.Range("D5").Value = 1.070546143554303
$ws.Range("E5").Value = 1.075854451335134
$ws.Range("F5").Value = 1.084471235129798
$ws.Range("I5").Value = 1.048632607461891
$ws.Range("J5").Value = 1.076794357895216
$ws.Range("K5").Value = 1.072895834606119
$ws.Range("L5").Value = 1.078191969401763
$ws.Range("M5").Value = 1.08678926481779

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.072595937746634
$ws.Range("D6").Value = 1.070601204865636
$ws.Range("E6").Value = 1.075917515478835
$ws.Range("F6").Value = 1.084535946167165
$ws.Range("I6").Value = 1.048650468915691
$ws.Range("J6").Value = 1.07685025026789
$ws.Range("K6").Value = 1.072942373128134
$ws.Range("L6").Value = 1.078246539356466
$ws.Range("M6").Value = 1.0868455531234

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.072102550357702
$ws.Range("D7").Value = 1.070222471415708
$ws.Range("E7").Value = 1.075483784253471
$ws.Range("F7").Value = 1.084090890245865
$ws.Range("I7").Value = 1.04852746046719
$ws.Range("J7").Value = 1.076465776254161
$ws.Range("K7").Value = 1.072622197402695
$ws.Range("L7").Value = 1.077871171833164
$ws.Range("M7").Value = 1.086458369879835

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.070037401668493
$ws.Range("D8").Value = 1.068636364429712
$ws.Range("E8").Value = 1.073668566370885
$ws.Range("F8").Value = 1.082228328075392
$ws.Range("I8").Value = 1.048008563015973
$ws.Range("J8").Value = 1.07485502088935
$ws.Range("K8").Value = 1.071279693808313
$ws.Range("L8").Value = 1.076298813183482
$ws.Range("M8").Value = 1.084836622864854

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.066386984248655
$ws.Range("D9").Value = 1.065829538792001
$ws.Range("E9").Value = 1.070460858742751
$ws.Range("F9").Value = 1.078937156792816
$ws.Range("I9").Value = 1.047076170300368
$ws.Range("J9").Value = 1.072002246235454
$ws.Range("K9").Value = 1.068897783719136
$ws.Range("L9").Value = 1.073514972110979
$ws.Range("M9").Value = 1.081965748771482

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.063945221332589
$ws.Range("D10").Value = 1.063949981410958
$ws.Range("E10").Value = 1.068315896080534
$ws.Range("F10").Value = 1.076736523396183
$ws.Range("I10").Value = 1.04644230763377
$ws.Range("J10").Value = 1.070090317454965
$ws.Range("K10").Value = 1.067298602830459
$ws.Range("L10").Value = 1.071649872603068
$ws.Range("M10").Value = 1.080042621920126

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.062885863023801
$ws.Range("D11").Value = 1.063134058149238
$ws.Range("E11").Value = 1.067385474704331
$ws.Range("F11").Value = 1.075781988436041
$ws.Range("I11").Value = 1.046164895434651
$ws.Range("J11").Value = 1.069259954167601
$ws.Range("K11").Value = 1.066603402380071
$ws.Range("L11").Value = 1.070839997973872
$ws.Range("M11").Value = 1.079207618948811

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.062492050905421
$ws.Range("D12").Value = 1.062830671409696
$ws.Range("E12").Value = 1.067039620850492
$ws.Range("F12").Value = 1.075427176476295
$ws.Range("I12").Value = 1.046061406748231
$ws.Range("J12").Value = 1.068951139315854
$ws.Range("K12").Value = 1.066344755399513
$ws.Range("L12").Value = 1.070538825601944
$ws.Range("M12").Value = 1.078897112415614

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.062576539550146
$ws.Range("D13").Value = 1.062895763325252
$ws.Range("E13").Value = 1.067113819292481
$ws.Range("F13").Value = 1.075503296520107
$ws.Range("I13").Value = 1.046083625627198
$ws.Range("J13").Value = 1.069017398510498
$ws.Range("K13").Value = 1.066400255097176
$ws.Range("L13").Value = 1.070603444001568
$ws.Range("M13").Value = 1.078963733042488

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.062853316962211
$ws.Range("D14").Value = 1.063108986602065
$ws.Range("E14").Value = 1.067356891542147
$ws.Range("F14").Value = 1.075752664819922
$ws.Range("I14").Value = 1.046156350131451
$ws.Range("J14").Value = 1.069234435250793
$ws.Range("K14").Value = 1.066582031115146
$ws.Range("L14").Value = 1.070815110130437
$ws.Range("M14").Value = 1.07918195954871

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.063023806054483
$ws.Range("D15").Value = 1.063240318359678
$ws.Range("E15").Value = 1.067506622495925
$ws.Range("F15").Value = 1.075906274810389
$ws.Range("I15").Value = 1.046201098985305
$ws.Range("J15").Value = 1.069368108042911
$ws.Range("K15").Value = 1.066693973665793
$ws.Range("L15").Value = 1.070945478172695
$ws.Range("M15").Value = 1.079316369604361

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.06401548322899
$ws.Range("D16").Value = 1.064004087530282
$ws.Range("E16").Value = 1.068377609901527
$ws.Range("F16").Value = 1.076799837385705
$ws.Range("I16").Value = 1.046460656261398
$ws.Range("J16").Value = 1.070145372923018
$ws.Range("K16").Value = 1.067344682638507
$ws.Range("L16").Value = 1.071703572786452
$ws.Range("M16").Value = 1.080097989739288

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.064636977903013
$ws.Range("D17").Value = 1.064482622888016
$ws.Range("E17").Value = 1.068923513174422
$ws.Range("F17").Value = 1.077359899327206
$ws.Range("I17").Value = 1.046622679042363
$ws.Range("J17").Value = 1.070632259905131
$ws.Range("K17").Value = 1.067752115997849
$ws.Range("L17").Value = 1.0721784915687
$ws.Range("M17").Value = 1.080587664851915

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.064999287258952
$ws.Range("D18").Value = 1.064761545915783
$ws.Range("E18").Value = 1.069241771766001
$ws.Range("F18").Value = 1.077686415684521
$ws.Range("I18").Value = 1.046716900192014
$ws.Range("J18").Value = 1.070916013541505
$ws.Range("K18").Value = 1.067989500596425
$ws.Range("L18").Value = 1.072455284718399
$ws.Range("M18").Value = 1.080873064919042

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.065122792106597
$ws.Range("D19").Value = 1.064856617982439
$ws.Range("E19").Value = 1.069350263231813
$ws.Range("F19").Value = 1.077797722710774
$ws.Range("I19").Value = 1.046748979085995
$ws.Range("J19").Value = 1.071012725835527
$ws.Range("K19").Value = 1.068070397953061
$ws.Range("L19").Value = 1.072549627109885
$ws.Range("M19").Value = 1.080970342014645

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.064570317955611
$ws.Range("D20").Value = 1.064431301158079
$ws.Range("E20").Value = 1.068864959272299
$ws.Range("F20").Value = 1.077299826411671
$ws.Range("I20").Value = 1.046605324927663
$ws.Range("J20").Value = 1.070580046346627
$ws.Range("K20").Value = 1.067708429646278
$ws.Range("L20").Value = 1.072127559976213
$ws.Range("M20").Value = 1.080535150065981

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.062771821806728
$ws.Range("D21").Value = 1.06304620645238
$ws.Range("E21").Value = 1.067285319895461
$ws.Range("F21").Value = 1.07567923917845
$ws.Range("I21").Value = 1.046134946897569
$ws.Range("J21").Value = 1.069170533917316
$ws.Range("K21").Value = 1.066528514213712
$ws.Range("L21").Value = 1.070752789419867
$ws.Range("M21").Value = 1.079117706981381

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.061639183684738
$ws.Range("D22").Value = 1.062173507868744
$ws.Range("E22").Value = 1.066290664878292
$ws.Range("F22").Value = 1.074658831363245
$ws.Range("I22").Value = 1.04583662289388
$ws.Range("J22").Value = 1.06828210964512
$ws.Range("K22").Value = 1.065784229649726
$ws.Range("L22").Value = 1.069886394517042
$ws.Range("M22").Value = 1.078224479959747

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.06223979549454
$ws.Range("D23").Value = 1.062636318027658
$ws.Range("E23").Value = 1.066818092562363
$ws.Range("F23").Value = 1.07519991180682
$ws.Range("I23").Value = 1.045995015501888
$ws.Range("J23").Value = 1.06875329197244
$ws.Range("K23").Value = 1.066179020928144
$ws.Range("L23").Value = 1.070345880935568
$ws.Range("M23").Value = 1.07869819082114

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.064600439328428
$ws.Range("D24").Value = 1.064454491852482
$ws.Range("E24").Value = 1.068891417747321
$ws.Range("F24").Value = 1.077326971265789
$ws.Range("I24").Value = 1.046613167382739
$ws.Range("J24").Value = 1.07060364014585
$ws.Range("K24").Value = 1.067728170445238
$ws.Range("L24").Value = 1.07215057444726
$ws.Range("M24").Value = 1.080558879914635

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.067332103586538
$ws.Range("D25").Value = 1.066556615598235
$ws.Range("E25").Value = 1.071291242741622
$ws.Range("F25").Value = 1.079789122874741
$ws.Range("I25").Value = 1.047319367627269
$ws.Range("J25").Value = 1.072741504936169
$ws.Range("K25").Value = 1.069515522435141
$ws.Range("L25").Value = 1.074236257133368
$ws.Range("M25").Value = 1.082709534673981
